# "merubah konfigurasi import data absensi"
#
# The attendance rows for Person ID 1 (rows 2-101, column B) were imported
# under the wrong account ("Nur A"); re-point them at "Admin" instead. The
# other two blocks of rows (Person ID 2 -> "Approve 01", Person ID 3 ->
# "User 01") keep their text - only this first block's label changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B101").Value = "Admin"

# Reset the view: scroll back up to the top of the sheet and select B3.
$win = $excel.ActiveWindow
try { $win.ScrollRow = 2 } catch {}
try { $win.ScrollColumn = 1 } catch {}
[void]$ws.Range("B3").Select()
